# Updated result of RQ1 for Ctags.
$wb = $excel.ActiveWorkbook

# --- RQ1 sheet: new Ctags numbers (row 3) and new Brlcad numbers (row 4) ---
$ws1 = $wb.Worksheets.Item("RQ1")

# Row 3 (Ctags): denominators/second-metric numerator change, numerators (C3) unchanged
$ws1.Range("D3").Value = 39
$ws1.Range("F3").Value = 208
$ws1.Range("G3").Value = 39

# Row 4 (Brlcad): new counts
$ws1.Range("C4").Value = 63
$ws1.Range("D4").Value = 26
$ws1.Range("F4").Value = 113
$ws1.Range("G4").Value = 26

# Row 7 (Jabref): data no longer available -> blank out, formulas fall back to #DIV/0!
$ws1.Range("C7:D7").ClearContents()
$ws1.Range("F7:G7").ClearContents()

# --- RQ2 / RQ3 / RQ4: clear the now-unavailable source counts, keep the formulas ---
$ws2 = $wb.Worksheets.Item("RQ2")
$ws2.Range("B3:C3").ClearContents()
$ws2.Range("E3:F3").ClearContents()
$ws2.Range("B4:C4").ClearContents()
$ws2.Range("E4:F4").ClearContents()
[void]$ws2.Range("A9").Select()

$ws3 = $wb.Worksheets.Item("RQ3")
$ws3.Range("B3:C3").ClearContents()
$ws3.Range("E3:F3").ClearContents()
$ws3.Range("B4:C4").ClearContents()
$ws3.Range("E4:F4").ClearContents()
[void]$ws3.Range("A9").Select()

$ws4 = $wb.Worksheets.Item("RQ4")
$ws4.Range("B3:C3").ClearContents()
$ws4.Range("E3:F3").ClearContents()
$ws4.Range("B4:C4").ClearContents()
$ws4.Range("E4:F4").ClearContents()
[void]$ws4.Range("A9").Select()

# --- Remove the now-redundant RQ1a sheet (its content was folded into RQ1) ---
$excel.DisplayAlerts = $false
$ws1a = $wb.Worksheets.Item("RQ1a")
[void]$ws1a.Delete()

# --- Make RQ1 the active/selected sheet, with A9 selected ---
[void]$ws1.Activate()
[void]$ws1.Range("A9").Select()
